$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")
$ws.Range("A1").Value = "TEST"
